$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 10-12 summary updates ----
$ws.Range("A10").Value = "No."
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 18
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 28
$ws.Range("A11").Value = "Marking"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("A12").Value = "Total"
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 72
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "70/112"

# ---- Rows 16-40: fill "Student Ans" column A (and D where applicable) ----
$ws.Range("A16").Value = "Option A"
$ws.Range("A16").Style = "correctStyle"
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"
$ws.Range("A17").Value = "Option D"
$ws.Range("A17").Style = "correctStyle"
$ws.Range("D17").Value = "Option C"
$ws.Range("D17").Style = "correctStyle"
$ws.Range("A18").Value = "Option B"
$ws.Range("A18").Style = "correctStyle"
$ws.Range("A19").Value = "Option C"
$ws.Range("A19").Style = "correctStyle"
$ws.Range("A21").Value = "Option C"
$ws.Range("A21").Style = "correctStyle"
$ws.Range("A22").Value = "Option D"
$ws.Range("A22").Style = "correctStyle"
$ws.Range("A24").Value = "Option A"
$ws.Range("A24").Style = "correctStyle"
$ws.Range("A27").Value = "Option A"
$ws.Range("A27").Style = "correctStyle"
$ws.Range("A28").Value = "Option D"
$ws.Range("A28").Style = "correctStyle"
$ws.Range("A29").Value = "Option D"
$ws.Range("A29").Style = "correctStyle"
$ws.Range("A30").Value = "Option B"
$ws.Range("A30").Style = "correctStyle"
$ws.Range("A31").Value = "Option D"
$ws.Range("A31").Style = "correctStyle"
$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"
$ws.Range("A33").Value = "Option C"
$ws.Range("A33").Style = "incorrectStyle"
$ws.Range("A35").Value = "Option B"
$ws.Range("A35").Style = "incorrectStyle"
$ws.Range("A38").Value = "Option A"
$ws.Range("A38").Style = "correctStyle"
$ws.Range("A39").Value = "Option D"
$ws.Range("A39").Style = "correctStyle"
$ws.Range("A40").Value = "Option D"
$ws.Range("A40").Style = "correctStyle"

# ---- Clear removed cells (columns G:H entirely, and D:E tail) ----
$ws.Range("G15:H21").Clear()
$ws.Range("D19:E40").Clear()
